$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 106, shifting existing rows 106-227 down to 107-228
$ws.Rows.Item(106).Insert()

# Populate the new row 106 with the new data record
$ws.Range("A106").Value = 11
$ws.Range("B106").Value = "Vega Monumental Concepción"
$ws.Range("C106").Value = "Bíobío"
$ws.Range("D106").Value = "2022-03-18"
$ws.Range("E106").Value = 8
$ws.Range("F106").Value = 100112017
$ws.Range("G106").Value = "Apio"
$ws.Range("H106").Value = "Americana (o)"
$ws.Range("I106").Value = "Primera"
$ws.Range("J106").Value = 300
$ws.Range("K106").Value = 6500
$ws.Range("L106").Value = 7000
$ws.Range("M106").Value = 6750
$ws.Range("N106").Value = "`$/docena de matas"
$ws.Range("O106").Value = "Región de Coquimbo"
$ws.Range("P106").Value = 1125
$ws.Range("Q106").Value = 6
$ws.Range("R106").Value = "Hortaliza"


